$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 was an empty template row; fill it in with a new work-day entry,
# matching the date-cell style used by the rows above it (A2:A5).
$ws.Range("A5").Copy($ws.Range("A6"))
$ws.Range("A6").Value = (Get-Date -Year 2020 -Month 6 -Day 2 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("B6").Value = "8 hours 30 minuts"

# Move the active selection to B7, as in the target workbook
$ws.Range("B7").Select()
